$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.624.45"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").Value = "3.104.62"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'584.96"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "'145.17"
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.100.12"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +7.49%  "

$ws.Range("D11").Value = "'5.65"
$ws.Range("E11").Value = "  -2.41%  "

$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -2.01%  "

$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").Value = "'36.94"
$ws.Range("E14").Value = "  +4.02%  "

$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "3.617.82"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").Value = "63.429.49"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").Value = "'7.07"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").Value = "3.096.73"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").Value = "'461.39"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "'0.723"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("D23").Value = "'7.43"
$ws.Range("E23").Value = "  -1.15%  "

$ws.Range("D24").Value = "'12.93"
$ws.Range("E24").Value = "  -2.73%  "

$ws.Range("D25").Value = "'81.19"
$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  +3.15%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "'9.23"
$ws.Range("E28").Value = "  +10.37%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").Value = "'6.93"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").Value = "'26.71"
$ws.Range("E34").Value = "  -0.80%  "

$ws.Range("D35").Value = "0.0₃0856"
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'1.03"
$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.38"
$ws.Range("E37").Value = "  +2.16%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'2.31"
$ws.Range("E38").Value = "  -4.14%  "

$ws.Range("E39").Value = "  -0.49%  "

$ws.Range("D40").Value = "'50.28"
$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("D41").Value = "'435.17"
$ws.Range("E41").Value = "  +1.13%  "

$ws.Range("D42").Value = "'8.69"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D43").Value = "'0.0369"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").Value = "2.881.85"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").Value = "'0.274"
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("E46").Value = "  -2.59%  "

$ws.Range("D47").Value = "'36.20"
$ws.Range("E47").Value = "  +2.68%  "

$ws.Range("D48").Value = "'125.64"
$ws.Range("E48").Value = "  +2.23%  "

$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("D51").Value = "'24.08"
$ws.Range("E51").Value = "  -1.60%  "
